# input data name change
#
# "industry" sheet: a new aperc_fuel -> ipcc_fuel mapping row for
# "15_solid_biomass$15_04_black_liquor" (black liquor) was missing.
# Previously row 466 ("15_solid_biomass$15_02_bagasse") was incorrectly
# pointing at the "Sulphite Lyes (Black Liquor)" ipcc_fuel. Insert the
# missing black-liquor row right after it and fix row 466 to point at
# the correct ipcc_fuel ("Other Primary Solid Biomass").
#
# "transformation" sheet: append a missing row mapping CHP plants /
# fuelwood-and-woodwaste to "1.A.1 - Energy Industries" / "Wood/Wood Waste".

$wb = $excel.ActiveWorkbook

$wsIndustry = $wb.Worksheets.Item("industry")

# Shift rows 467.. down by one to make room for the new row.
$wsIndustry.Rows.Item(467).Insert()

# Fix the now-mismatched ipcc_fuel for the bagasse row (466).
$wsIndustry.Range("D466").Value = 'Other Primary Solid Biomass'

# Fill in the newly inserted row (467) with the black-liquor mapping.
$wsIndustry.Range("A467").Value = '14_industry_sector$14_03_manufacturing$14_03_02_chemical_incl_petrochemical$14_03_02_01_fs'
$wsIndustry.Range("B467").Value = '15_solid_biomass$15_04_black_liquor'
$wsIndustry.Range("C467").Value = '1.A.2 - Manufacturing Industries and Construction'
$wsIndustry.Range("D467").Value = 'Sulphite Lyes (Black Liquor)'

$wsTransformation = $wb.Worksheets.Item("transformation")

# Append the missing CHP-plants / fuelwood-and-woodwaste row at the end.
$wsTransformation.Range("A169").Value = '09_total_transformation_sector$09_02_chp_plants$x$x'
$wsTransformation.Range("B169").Value = '15_solid_biomass$15_01_fuelwood_and_woodwaste'
$wsTransformation.Range("C169").Value = '1.A.1 - Energy Industries'
$wsTransformation.Range("D169").Value = 'Wood/Wood Waste'
